# Add the new "Initials" / "Complete words" rows to the recap-rephrasing
# table (Table1) on Sheet1, extending it from A1:B71 to A1:B81.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$newEntries = @(
    @("L/D",   "lay down"),
    @("LMP",   "liquid mud plant"),
    @("HLB",   "halliburton"),
    @("CONT.", "continue"),
    @("PU",    "pick up"),
    @("MU",    "Make up"),
    @("NU",    "nipple up"),
    @("RU",    "rig up"),
    @("RD",    "rig down"),
    @("ND",    "nipple down")
)

foreach ($entry in $newEntries) {
    $row = $lo.ListRows.Add()
    $rowIndex = $row.Range.Row
    $ws.Cells.Item($rowIndex, 1).Value2 = $entry[0]
    $ws.Cells.Item($rowIndex, 2).Value2 = $entry[1]
}

# Match the saved view state: scrolled down with B83 (the first empty row
# in column B following the new data) selected.
[void]$ws.Range("A59").Select()
[void]$ws.Range("B83").Select()
